# "Generate Report for Archive"
#
# Localization status moved on from handoff: every cell that used to read
# "Ready for handoff" now reads "In Translation". That text lives in the
# "zh-cn"/"de-de" columns (E/F) of the Overview sheet, and in the "Status"
# column (C) of each per-language sheet. Because the new status text is
# shorter, those columns are narrower than before.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# New column width (character units) for the narrower Status columns.
# The runtime quantizes ColumnWidth to sixths of a character, so 12.5 is
# the input that lands on the closest achievable width to the original
# file's target width.
$newColumnWidth = 12.5

# Overview sheet: the "zh-cn" and "de-de" columns (E and F) hold the
# status text and shrink along with it.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth

# zh-cn sheet: column C is "Status".
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth

# de-de sheet: column C is "Status".
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
